# Add new row 22 ("DRA100") to the "Test Cases" sheet for the
# Forgot-password test scripts (DRA / IPA application).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 22 from the row above so it inherits the same banded
# border / wrap-text formatting already used throughout the table,
# then overwrite the copied values with the new test-case data.
$ws.Range("A21:E21").Copy($ws.Range("A22:E22"))
$ws.Range("B10").Copy($ws.Range("B22"))
$ws.Range("C20").Copy($ws.Range("C22"))

$ws.Range("A22").Value = "DRA100"

$ws.Range("B22").Value = @'
OPQA-1934||
OPQA-1935&OPQA-3687||
OPQA-4230||OPQA-4229||
OPQA-4231||OPQA-4232||
OPQA-4636||OPQA-4261||
OPQA-4244||OPQA-4264||
OPQA-4265||OPQA-4237||
OPQA-4239||OPQA-4240||
OPQA-4246||OPQA-4248||
OPQA-4252
'@

$ws.Range("C22").Value = @'
Verify that Forgot your password? Link is clickable on NEON Landing page and End note landing page||Verify that the system is navigating to Forgot Password page or not, after clicking on Forgot your password? Link&Verify that,the system should support a ENW password reset workflow with the following configurations||Verify that system should not inform user that entered email is not found.||Verify that user should be able to enter email address in Forgot password page.||Verify that  forget password service should send a forgot password email when the email entered is registered in the system||Verify that the platform password reset service should send a platform forget password email with branding that corresponds with the originating application as per wireframe||Verify that When the password reset token in the email is valid, upon clicking the password reset link in the the platform forget password email, the user shall be taken to the External Password Reset Page||Verify that External Password Reset Page should have a new password field where the user enters their new password.||Verify that when reset Password Token already used user should be taken to sign in screen||Verify that upon successful submission of a password change, The user should receive a password change confirmation email to the user's primary email address with branding that corresponds with the application that the user completed the password change||Verify that the password change confirmation email should reference the fact that credentials are shared across all products.||Verify that when the password reset token in the email is expired or already used, upon clicking the password reset link in the the platform forget password email, the user should be taken to the External Invalid Password Reset Token Page.||Verify that the email address on the External Invalid Password Reset Token Page should be pre-populated with the email address that matches the email that the forgot password email was sent.||Verify that user who clicks the submit button on the the External Invalid Password Reset Token page, should be taken to the target application sign in page.||Verify that when Email address is known from password reset token,error message 'The email address is prepopulated.' should be displayed and email address field should be editable||Verify that when Email address is not known from password reset token,email address field should be blank and user should be able to enter any email address||Verify that error message Please enter a valid email address.should be displayed in red color when user enters email address in wrong format
'@

$ws.Range("D22").Value = "Y"
# E22 (Results) is left blank, same as the other data rows.

# Match the auto-fitted row height Excel computes for the wrapped text.
$ws.Rows.Item(22).RowHeight = 172.8

$ws.Range("A22:E22").Select()
